# Applies the crypto-price refresh for Sun Jun 23 04:31:36 UTC 2024
# (rows 2-51 of Sheet1): updated Price/Volume(1h) values, plus the
# Stacks/ONDO row swap at rows 46-47.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.378.14'
$ws.Range('E2').Value = '  +0.19%  '

$ws.Range('D3').Value = '3.517.69'
$ws.Range('E3').Value = '  +0.41%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.70'
$ws.Range('E5').Value = '  +1.23%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.82'
$ws.Range('E6').Value = '  -0.02%  '

$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('E8').Value = '  +0.32%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.63'
$ws.Range('E9').Value = '  +6.81%  '

$ws.Range('E10').Value = '  +0.14%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.390'
$ws.Range('E11').Value = '  +3.95%  '

$ws.Range('D12').Value = '4.119.39'
$ws.Range('E12').Value = '  +0.51%  '

$ws.Range('E13').Value = '  +1.42%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000181'
$ws.Range('E14').Value = '  +1.01%  '

$ws.Range('D15').Value = '3.518.65'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.93'
$ws.Range('E16').Value = '  -1.74%  '

$ws.Range('D17').Value = '64.380.32'
$ws.Range('E17').Value = '  +0.18%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.00'
$ws.Range('E18').Value = '  +2.49%  '

$ws.Range('E19').Value = '  +3.45%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.59'
$ws.Range('E20').Value = '  -1.91%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '394.34'
$ws.Range('E21').Value = '  +2.74%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.577'
$ws.Range('E22').Value = '  +1.43%  '

$ws.Range('D23').Value = '3.660.09'
$ws.Range('E23').Value = '  +0.54%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.71'
$ws.Range('E24').Value = '  +1.04%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.998'
$ws.Range('E25').Value = '  -0.23%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.74'

$ws.Range('E27').Value = '  +3.14%  '

$ws.Range('E28').Value = '  -0.04%  '

$ws.Range('E29').Value = '  -1.65%  '

$ws.Range('E30').Value = '  +1.52%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.31'
$ws.Range('E31').Value = '  -0.09%  '

$ws.Range('E32').Value = '  -6.40%  '

$ws.Range('E33').Value = '  +7.67%  '

$ws.Range('D34').Value = '3.549.25'
$ws.Range('E34').Value = '  +0.76%  '

$ws.Range('E35').Value = '  +0.04%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.45'
$ws.Range('E36').Value = '  -0.77%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.36'
$ws.Range('E37').Value = '  +0.70%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.97'
$ws.Range('E38').Value = '  +1.63%  '

$ws.Range('E39').Value = '  +1.65%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '166.93'
$ws.Range('E40').Value = '  +1.49%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0792'
$ws.Range('E41').Value = '  +1.21%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.813'
$ws.Range('E42').Value = '  +0.38%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '25.52'
$ws.Range('E43').Value = '  -1.97%  '

$ws.Range('E44').Value = '  -0.01%  '

$ws.Range('E45').Value = '  +0.96%  '

$ws.Range('B46').Value = 'ONDO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.18'
$ws.Range('E46').Value = '  -1.85%  '

$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.66'
$ws.Range('E47').Value = '  +2.95%  '

$ws.Range('E48').Value = '  +0.67%  '

$ws.Range('D49').Value = '2.406.77'
$ws.Range('E49').Value = '  -2.63%  '

$ws.Range('E50').Value = '  -2.13%  '

$ws.Range('E51').Value = '  +0.34%  '
